$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.981.48'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '2.973.19'
$ws.Range('E3').Value = '  +4.35%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '353.93'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.68%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '112.93'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.28%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.567'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +2.11%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +1.96%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '39.77'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.16%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0895'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +5.12%  '
$ws.Range('E12').Value = '  +1.28%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '20.12'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.09%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '7.92'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').Value = '3.431.12'
$ws.Range('E15').Value = '  +4.20%  '
$ws.Range('D16').Value = '2.970.02'
$ws.Range('E16').Value = '  +4.21%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.992'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '52.055.41'
$ws.Range('E18').Value = '  +0.27%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.69'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.48%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '3.34'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.61%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '14.45'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +6.41%  '
$ws.Range('D22').Value = '0.0₃0991'
$ws.Range('E22').Value = '  +1.57%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '71.47'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.30%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '270.31'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.48%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.81'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.72%  '
$ws.Range('E26').Value = '  +10.31%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '27.28'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +3.55%  '
$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.115'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +28.40%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.07%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '7.52'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +19.32%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '10.77'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.96%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '37.70'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -3.87%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '2.29'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.91%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '6.22'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +10.32%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '53.08'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.68%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0453'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  -0.16%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.37'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.92%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '18.99'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.53%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.07'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +2.41%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.68'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +5.83%  '
$ws.Range('E42').Value = '  +1.68%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '23.73'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +6.60%  '
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.55'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.53'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.04%  '
$ws.Range('D47').Value = '2.183.88'
$ws.Range('E47').Value = '  +0.51%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '114.01'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -7.13%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.246'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.19%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0343'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +8.98%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.939'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.22%  '
